# Requirements loaded into dictionaries
# Insert a new header/title row at the very top of the sheet and push the
# existing course list down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above row 1 (everything below shifts down automatically,
# and the shared-string indices used by the existing cells are preserved).
$ws.Rows.Item(1).Insert()

# Populate the new first row with the new title text.
$ws.Range("A1").Value = "Cultural Diversity (CD) Courses Offered in 2016-2017"

# Give the new title cell an explicit black font color (as opposed to the
# theme-based color used by the rest of the sheet), which creates a second
# font / cell style entry, matching the authored change.
$ws.Range("A1").Font.Color = 0
